$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "B" and "C" sub-rows within each year block (columns A:E).
# (The "A" and "D" rows of each 4-row year block are unaffected.)
# rows 3 <-> 4
$ws.Cells.Item(3,1).Value = "2000年C"
$ws.Cells.Item(3,2).Value = 98.09999999999999
$ws.Cells.Item(3,3).ClearContents()
$ws.Cells.Item(3,4).Value = 50.7
$ws.Cells.Item(3,5).Value = 1831.3
$ws.Cells.Item(4,1).Value = "2000年B"
$ws.Cells.Item(4,2).Value = 97.5
$ws.Cells.Item(4,3).ClearContents()
$ws.Cells.Item(4,4).Value = 39.2
$ws.Cells.Item(4,5).Value = 1051

# rows 7 <-> 8
$ws.Cells.Item(7,1).Value = "2001年C"
$ws.Cells.Item(7,2).Value = 97.59999999999999
$ws.Cells.Item(7,3).Value = -0.5
$ws.Cells.Item(7,4).Value = 57.8
$ws.Cells.Item(7,5).Value = 1847.6
$ws.Cells.Item(8,1).Value = "2001年B"
$ws.Cells.Item(8,2).Value = 98.7
$ws.Cells.Item(8,3).ClearContents()
$ws.Cells.Item(8,4).Value = 19.2
$ws.Cells.Item(8,5).Value = 1109.7

# rows 11 <-> 12
$ws.Cells.Item(11,1).Value = "2002年C"
$ws.Cells.Item(11,2).Value = 99.7
$ws.Cells.Item(11,3).ClearContents()
$ws.Cells.Item(11,4).Value = 8.6
$ws.Cells.Item(11,5).Value = 1981.1
$ws.Cells.Item(12,1).Value = "2002年B"
$ws.Cells.Item(12,2).Value = 98.90000000000001
$ws.Cells.Item(12,3).ClearContents()
$ws.Cells.Item(12,4).Value = 12.7
$ws.Cells.Item(12,5).Value = 1169.7

# rows 15 <-> 16
$ws.Cells.Item(15,1).Value = "2003年C"
$ws.Cells.Item(15,2).Value = 99.40000000000001
$ws.Cells.Item(15,3).ClearContents()
$ws.Cells.Item(15,4).Value = 8.800000000000001
$ws.Cells.Item(15,5).Value = 2070.2
$ws.Cells.Item(16,1).Value = "2003年B"
$ws.Cells.Item(16,2).Value = 98.40000000000001
$ws.Cells.Item(16,3).Value = -0.5
$ws.Cells.Item(16,4).Value = 15.2
$ws.Cells.Item(16,5).Value = 1171.4

# rows 19 <-> 20
$ws.Cells.Item(19,1).Value = "2004年C"
$ws.Cells.Item(19,2).Value = 98.3
$ws.Cells.Item(19,3).Value = -1.2
$ws.Cells.Item(19,4).Value = 57.7
$ws.Cells.Item(19,5).Value = 2322
$ws.Cells.Item(20,1).Value = "2004年B"
$ws.Cells.Item(20,2).Value = 99.3
$ws.Cells.Item(20,3).Value = 0.9
$ws.Cells.Item(20,4).Value = 14.4
$ws.Cells.Item(20,5).Value = 1394.4

# rows 23 <-> 24
$ws.Cells.Item(23,1).Value = "2005年C"
$ws.Cells.Item(23,2).Value = 98.40000000000001
$ws.Cells.Item(23,3).Value = 0.1
$ws.Cells.Item(23,4).Value = 56
$ws.Cells.Item(23,5).Value = 2455.9
$ws.Cells.Item(24,1).Value = "2005年B"
$ws.Cells.Item(24,2).Value = 98.59999999999999
$ws.Cells.Item(24,3).Value = -0.7
$ws.Cells.Item(24,4).Value = 24.2
$ws.Cells.Item(24,5).Value = 1420.4

# rows 27 <-> 28
$ws.Cells.Item(27,1).Value = "2006年C"
$ws.Cells.Item(27,2).Value = 98.59999999999999
$ws.Cells.Item(27,3).Value = 0
$ws.Cells.Item(27,4).Value = 34.9
$ws.Cells.Item(27,5).Value = 2789.1
$ws.Cells.Item(28,1).Value = "2006年B"
$ws.Cells.Item(28,2).Value = 98.8
$ws.Cells.Item(28,3).Value = 0.2
$ws.Cells.Item(28,4).Value = 27.8
$ws.Cells.Item(28,5).Value = 1618

# rows 31 <-> 32
$ws.Cells.Item(31,1).Value = "2007年C"
$ws.Cells.Item(31,2).Value = 99.5
$ws.Cells.Item(31,3).Value = 1.1
$ws.Cells.Item(31,4).Value = 25.6
$ws.Cells.Item(31,5).Value = 3162.1
$ws.Cells.Item(32,1).Value = "2007年B"
$ws.Cells.Item(32,2).Value = 99.09999999999999
$ws.Cells.Item(32,3).Value = 0.3
$ws.Cells.Item(32,4).Value = 24.1
$ws.Cells.Item(32,5).Value = 1850.4

# rows 35 <-> 36
$ws.Cells.Item(35,1).Value = "2008年C"
$ws.Cells.Item(35,2).Value = 99
$ws.Cells.Item(35,3).Value = -0.2
$ws.Cells.Item(35,4).Value = 28.5
$ws.Cells.Item(35,5).Value = 3294.9
$ws.Cells.Item(36,1).Value = "2008年B"
$ws.Cells.Item(36,2).Value = 97.8
$ws.Cells.Item(36,3).Value = -0.3
$ws.Cells.Item(36,4).Value = 35.3
$ws.Cells.Item(36,5).Value = 1917.8

# rows 39 <-> 40
$ws.Cells.Item(39,1).Value = "2009年C"
$ws.Cells.Item(39,2).Value = 99.7
$ws.Cells.Item(39,3).Value = 0.2
$ws.Cells.Item(39,4).Value = 11.5
$ws.Cells.Item(39,5).Value = 3456.1
$ws.Cells.Item(40,1).Value = "2009年B"
$ws.Cells.Item(40,2).Value = 99
$ws.Cells.Item(40,3).Value = 1.8
$ws.Cells.Item(40,4).Value = 8.5
$ws.Cells.Item(40,5).Value = 2025.2

# rows 43 <-> 44
$ws.Cells.Item(43,1).Value = "2010年C"
$ws.Cells.Item(43,2).Value = 99.5
$ws.Cells.Item(43,3).Value = -0.6
$ws.Cells.Item(43,4).Value = 27.1
$ws.Cells.Item(43,5).Value = 3626.3
$ws.Cells.Item(44,1).Value = "2010年B"
$ws.Cells.Item(44,2).Value = 98.7
$ws.Cells.Item(44,3).Value = 0
$ws.Cells.Item(44,4).Value = 17
$ws.Cells.Item(44,5).Value = 2086

# rows 47 <-> 48
$ws.Cells.Item(47,1).Value = "2011年C"
$ws.Cells.Item(47,2).Value = 99.59999999999999
$ws.Cells.Item(47,3).Value = 0.6
$ws.Cells.Item(47,4).Value = 3.7
$ws.Cells.Item(47,5).Value = 3931.1
$ws.Cells.Item(48,1).Value = "2011年B"
$ws.Cells.Item(48,2).Value = 99.3
$ws.Cells.Item(48,3).Value = 0
$ws.Cells.Item(48,4).Value = 7.3
$ws.Cells.Item(48,5).Value = 2329

# rows 51 <-> 52
$ws.Cells.Item(51,1).Value = "2012年C"
$ws.Cells.Item(51,2).Value = 99.90000000000001
$ws.Cells.Item(51,3).Value = 1
$ws.Cells.Item(51,4).Value = 2.1
$ws.Cells.Item(51,5).Value = 3976.3
$ws.Cells.Item(52,1).Value = "2012年B"
$ws.Cells.Item(52,2).Value = 101.7
$ws.Cells.Item(52,3).Value = 2.1
$ws.Cells.Item(52,4).Value = -6.4
$ws.Cells.Item(52,5).Value = 2434.4

# rows 55 <-> 56
$ws.Cells.Item(55,1).Value = "2013年C"
$ws.Cells.Item(55,2).Value = 100.2
$ws.Cells.Item(55,3).Value = 0.2
$ws.Cells.Item(55,4).Value = 2
$ws.Cells.Item(55,5).Value = 4090.8
$ws.Cells.Item(56,1).Value = "2013年B"
$ws.Cells.Item(56,2).Value = 99.8
$ws.Cells.Item(56,3).Value = -0.5
$ws.Cells.Item(56,4).Value = 4.6
$ws.Cells.Item(56,5).Value = 2483.7

# rows 59 <-> 60
$ws.Cells.Item(59,1).Value = "2014年C"
$ws.Cells.Item(59,2).Value = 101.7
$ws.Cells.Item(59,3).Value = 1.9
$ws.Cells.Item(59,4).Value = -1.5
$ws.Cells.Item(59,5).Value = 4151.5
$ws.Cells.Item(60,1).Value = "2014年B"
$ws.Cells.Item(60,2).Value = 102
$ws.Cells.Item(60,3).Value = 2
$ws.Cells.Item(60,4).Value = -3.1
$ws.Cells.Item(60,5).Value = 2624.7

# rows 63 <-> 64
$ws.Cells.Item(63,1).Value = "2015年C"
$ws.Cells.Item(63,2).Value = 98.8
$ws.Cells.Item(63,3).Value = -2.1
$ws.Cells.Item(63,4).Value = 11.5
$ws.Cells.Item(63,5).Value = 3832.4
$ws.Cells.Item(64,1).Value = "2015年B"
$ws.Cells.Item(64,2).Value = 98.90000000000001
$ws.Cells.Item(64,3).Value = -2.5
$ws.Cells.Item(64,4).Value = 12.1
$ws.Cells.Item(64,5).Value = 2422.8

# Remove the now-redundant F (产销率) and G (销售量) columns
$ws.Columns("F:G").Delete()
